# Sprint 2 plan update
# ---------------------
# This reassigns work in the sprint backlog grid on Sheet1:
#   - Row 4 (Task T2): the Day 1/Day 2 owner tag changes from "A:1" (Alex) to "S:1" (Saad)
#   - Row 7 (Task T7): Story Points goes from 3 to 5, and it now also picks up "A:1"
#     (Alex) on Day 1 and Day 2, in addition to the Day 3-5 "A:1" it already had,
#     i.e. Alex now works the whole week on this task.
#   - Row 8 (Task T8): Story Points goes from 5 to 3, and it now loses "S:1" (Saad)
#     on Day 1 and Day 2 (those cells become blank), keeping "S:1" only on Day 3-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4: T2 owner tag Day1/Day2 (A:1 -> S:1) ---
$ws.Range("E4").Value = "S:1"
$ws.Range("F4").Value = "S:1"

# --- Row 7: T7 story points + Day1/Day2 owner tag (blank -> A:1) ---
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = "A:1"
$ws.Range("F7").Value = "A:1"

# --- Row 8: T8 story points + Day1/Day2 owner tag (S:1 -> blank) ---
$ws.Range("D8").Value = 3
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()

# Re-affirm the (unchanged) alignment already used by the A-column task-owner
# cells, so formatting stays exactly as intended after the edits above.
$ws.Range("A5").HorizontalAlignment = 1
$ws.Range("A5").VerticalAlignment = -4108

$ws.Range("A7:A10").HorizontalAlignment = -4131
$ws.Range("A7:A10").VerticalAlignment = -4108

$ws.Range("A11:A12").HorizontalAlignment = 1
$ws.Range("A11:A12").VerticalAlignment = -4108

# --- Update the sheet's last-used selection to reflect where the user ended
#     up working (E8:F8, with E8 active), matching the saved view state.
$ws.Activate()
$ws.Range("E8:F8").Select()
